$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds cells (rows 2-7) per diff ---
$ws.Range("F2").Value = 2.6
$ws.Range("G2").Value = 3
$ws.Range("J2").Value = 3.6
$ws.Range("K2").Value = 4.1
$ws.Range("L2").Value = 1.3
$ws.Range("N2").Value = 4.2
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 1.73
$ws.Range("R2").Value = 1.44
$ws.Range("S2").Value = 2.84
$ws.Range("U2").Value = 2.32
$ws.Range("V2").Value = 1.55
$ws.Range("Y2").Value = 14
$ws.Range("Z2").Value = 20
$ws.Range("AA2").Value = 44
$ws.Range("AB2").Value = 14.5
$ws.Range("AC2").Value = 9
$ws.Range("AF2").Value = 22
$ws.Range("AG2").Value = 13.5
$ws.Range("AH2").Value = 16.5
$ws.Range("AI2").Value = 42
$ws.Range("AJ2").Value = 50
$ws.Range("AK2").Value = 30
$ws.Range("AL2").Value = 44
$ws.Range("AM2").Value = 75
$ws.Range("AN2").Value = 26
$ws.Range("AO2").Value = 20
$ws.Range("T3").Value = 1.79
$ws.Range("Q4").Value = 1.7
$ws.Range("S4").Value = 2.72
$ws.Range("T4").Value = 1.59
$ws.Range("U4").Value = 2.36
$ws.Range("W4").Value = 1.65
$ws.Range("H5").Value = 1.76
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1000
$ws.Range("P7").Value = 1.91
$ws.Range("T7").Value = 1.71

# --- Add new row 8 (Colombian Primera A match) ---
$ws.Range("A8").Value = "Colombian Primera A"
# Leading apostrophe forces text so Excel doesn't coerce this into a date serial
$ws.Range("B8").Value = "'2025-12-16"
$ws.Range("C8").Value = "21:30:00"
$ws.Range("D8").Value = "Tolima"
$ws.Range("E8").Value = "Junior FC Barranquilla"
$ws.Range("F8").Value = 1.87
$ws.Range("G8").Value = 1.95
$ws.Range("H8").Value = 4.3
$ws.Range("I8").Value = 4.9
$ws.Range("J8").Value = 3.85
$ws.Range("K8").Value = 4
$ws.Range("L8").Value = 1.01
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 3.1
$ws.Range("O8").Value = 1.38
$ws.Range("P8").Value = 1.73
$ws.Range("Q8").Value = 2.1
$ws.Range("R8").Value = 1.27
$ws.Range("S8").Value = 3.9
$ws.Range("T8").Value = 1.94
$ws.Range("U8").Value = 1.87
$ws.Range("V8").Value = 1.26
$ws.Range("W8").Value = 2.04
$ws.Range("X8").Value = 15
$ws.Range("Y8").Value = 15
$ws.Range("Z8").Value = 36
$ws.Range("AA8").Value = 1000
$ws.Range("AB8").Value = 8
$ws.Range("AC8").Value = 8.800000000000001
$ws.Range("AD8").Value = 19.5
$ws.Range("AE8").Value = 75
$ws.Range("AF8").Value = 12
$ws.Range("AG8").Value = 11
$ws.Range("AH8").Value = 23
$ws.Range("AI8").Value = 1000
$ws.Range("AJ8").Value = 24
$ws.Range("AK8").Value = 24
$ws.Range("AL8").Value = 46
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 18
$ws.Range("AO8").Value = 120
